# Apply edits described by the diff: strip stale w:proofErr spell-check
# markers, merge split runs back into single runs, rewrite the "Concept"
# paragraph content, and relocate the _GoBack bookmark.
$d = $word.ActiveDocument

# Paragraph 2: "Indeling" - drop proofErr wrapper
$p2 = $d.Paragraphs(2).Range
$p2.InsertXML('<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:t>Indeling</w:t></w:r></w:p>')

# Paragraph 3: "Het idee..." - merge the asylum-related runs, drop proofErr
$p3 = $d.Paragraphs(3).Range
$p3.InsertXML('<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr><w:rPr><w:lang w:val="nl-NL"/></w:rPr></w:pPr><w:r w:rsidRPr="00E90E9C"><w:rPr><w:lang w:val="nl-NL"/></w:rPr><w:t xml:space="preserve">Het idee: je bent op de vlucht voor de mens, je bent aangekomen bij de safe haven die je zocht. </w:t></w:r><w:r><w:rPr><w:lang w:val="nl-NL"/></w:rPr><w:t>Om binnen te komen heb je 7 specifieke objecten nodig, deze krijg je door voor elk object een raadsel op te lossen. Als je de raadsels hebt opgelost en de juiste objecten hebt verzameld met daarmee bewezen hebt dat je weet waar je naartoe gaat krijg je asylum.</w:t></w:r></w:p>')

# Paragraph 4: "Concept..." - reworded content, drop proofErr, add bookmark
$p4 = $d.Paragraphs(4).Range
$p4.InsertXML('<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr><w:rPr><w:lang w:val="nl-NL"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="nl-NL"/></w:rPr><w:t>Concept: het concept is een first person scavenge hunt adventure puzzle game. Dit is omdat puur scavenge hunt of puur adventure het idee niet zou omschrijven. Mond vol maar het klopt. Je begint in de slums nadat je ben geintroduceerd in het spel. Hier weet je al dat je de raadsels moet opl</w:t></w:r><w:r><w:rPr><w:lang w:val="nl-NL"/></w:rPr><w:t>ossen om</w:t></w:r><w:r><w:rPr><w:lang w:val="nl-NL"/></w:rPr><w:t xml:space="preserve"> 7 veschillende objecten in 7 verschillende gebieden</w:t></w:r><w:r><w:rPr><w:lang w:val="nl-NL"/></w:rPr><w:t xml:space="preserve"> te vinden </w:t></w:r><w:r><w:rPr><w:lang w:val="nl-NL"/></w:rPr><w:t xml:space="preserve">, dit weet je </w:t></w:r><w:r><w:rPr><w:lang w:val="nl-NL"/></w:rPr><w:t>do</w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/><w:r><w:rPr><w:lang w:val="nl-NL"/></w:rPr><w:t>ordat</w:t></w:r><w:r><w:rPr><w:lang w:val="nl-NL"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:r><w:rPr><w:lang w:val="nl-NL"/></w:rPr><w:t>Malcom dit je verteld aan het begin</w:t></w:r><w:r><w:rPr><w:lang w:val="nl-NL"/></w:rPr><w:t xml:space="preserve">. </w:t></w:r><w:r><w:rPr><w:lang w:val="nl-NL"/></w:rPr><w:t xml:space="preserve">Eenmaal klaar hiermee ga je naar </w:t></w:r><w:r><w:rPr><w:lang w:val="nl-NL"/></w:rPr><w:t>Malcom</w:t></w:r><w:r><w:rPr><w:lang w:val="nl-NL"/></w:rPr><w:t xml:space="preserve"> waar je de objecten aan overhandigt om de city in de komen. </w:t></w:r><w:r><w:rPr><w:lang w:val="nl-NL"/></w:rPr><w:t>/</w:t></w:r><w:r><w:rPr><w:lang w:val="nl-NL"/></w:rPr><w:t>Binnen krijg je een aantal keuzes als welk huis en wat je rol is</w:t></w:r><w:r><w:rPr><w:lang w:val="nl-NL"/></w:rPr><w:t>/</w:t></w:r><w:r><w:rPr><w:lang w:val="nl-NL"/></w:rPr><w:t xml:space="preserve">, hiermee heb je je doel gehaald en is het spel klaar. </w:t></w:r></w:p>')

# Paragraph 6: "Semi-realistisch..." - merge into a single run, drop proofErr
# and remove the _GoBack bookmark (it moved to paragraph 4)
$p6 = $d.Paragraphs(6).Range
$p6.InsertXML('<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="ListParagraph"/><w:rPr><w:lang w:val="nl-NL"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="nl-NL"/></w:rPr><w:t xml:space="preserve">Semi-realistisch cartoony, omdat we een stap hoger moesten dan alleen cartoony. De characters zijn daarom nog steeds cartoony maar de omgeving zal semi-realistisch zijn. Verder hebben de slums een echt krottenwijk idee als je zou zien als je in een brazilië bent of als je door dying light speelt. </w:t></w:r></w:p>')

# Paragraph 8: "Ieder voor zich..." - merge into a single run, drop proofErr
$p8 = $d.Paragraphs(8).Range
$p8.InsertXML('<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="ListParagraph"/><w:rPr><w:lang w:val="nl-NL"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="nl-NL"/></w:rPr><w:t>Ieder voor zich in specifieke roles, producer doet de samenvatting.</w:t></w:r></w:p>')

Write-Output "applied edits"
